$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 599, shifting existing rows 599:613 down to 600:614
$ws.Rows.Item(599).Insert()

# Populate the newly inserted row 599 with the new weekly record.
# (Columns A,B,C,E,F,G,H,I,R are identical for every record in this block.)
$ws.Cells.Item(599, 1).Value = 8
$ws.Cells.Item(599, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(599, 3).Value = "Coquimbo"
$ws.Cells.Item(599, 4).Value = 45239
$ws.Cells.Item(599, 5).Value = 4
$ws.Cells.Item(599, 6).Value = 100112032
$ws.Cells.Item(599, 7).Value = "Zapallo italiano"
$ws.Cells.Item(599, 8).Value = "Sin especificar"
$ws.Cells.Item(599, 9).Value = "Primera"
$ws.Cells.Item(599, 10).Value = 500
$ws.Cells.Item(599, 11).Value = 9000
$ws.Cells.Item(599, 12).Value = 10000
$ws.Cells.Item(599, 13).Value = 9500
$ws.Cells.Item(599, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(599, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(599, 16).Value = 158
$ws.Cells.Item(599, 17).Value = 60
$ws.Cells.Item(599, 18).Value = "Hortaliza"

# Match the date-time number format used by the rest of column D.
$ws.Cells.Item(599, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
